$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45185
}
